$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "date" columns (B and E) actually hold plain text strings like
# "2020-09-06" in this workbook (t="inlineStr"), not real Excel dates.
# Force text formatting before writing so COM doesn't coerce the string
# into a date serial number.
$dateCells = "B3","E3","B4","E4","B5","E5","B6","E6","B7","E7","B8","E8"
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 3: A200001 - same client/service, quantities & totals changed
$ws.Range("B3").Value = "2020-09-12"
$ws.Range("E3").Value = "2020-09-12"
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 112

# Row 4: A200002 - client changed from "Luis Zurita Herrera"/DNI to "Clientes contado "
$ws.Range("B4").Value = "2020-09-12"
$ws.Range("C4").Value = "Clientes contado "
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "2020-09-12"
$ws.Range("G4").Value = 28
$ws.Range("J4").Value = 28

# Row 5: A200003 - client changed from "Luis Zurita Herrera"/DNI to "Clientes contado "
$ws.Range("B5").Value = "2020-09-12"
$ws.Range("C5").Value = "Clientes contado "
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "2020-09-12"
$ws.Range("G5").Value = 84
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 84

# Row 6: A200004 - only date/amount changed
$ws.Range("B6").Value = "2020-09-12"
$ws.Range("E6").Value = "2020-09-12"
$ws.Range("G6").Value = 28
$ws.Range("J6").Value = 28

# Row 7: A200005 - client changed from "Clientes contado " to "luis aslñdkf"/DNI
$ws.Range("B7").Value = "2020-09-12"
$ws.Range("C7").Value = "luis aslñdkf"
$ws.Range("D7").Value = "74666101M"
$ws.Range("E7").Value = "2020-09-12"
$ws.Range("G7").Value = 28
$ws.Range("J7").Value = 28

# Row 8: A200006 - only date/amount changed
$ws.Range("B8").Value = "2020-09-12"
$ws.Range("E8").Value = "2020-09-12"
$ws.Range("G8").Value = 28
$ws.Range("J8").Value = 28

# Rows 9-13 (A200007 .. A200011) are removed entirely; dimension shrinks to A1:J8
$ws.Range("A9:J13").Delete()
